$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (fills) from the last existing data row down to the new rows
$ws.Range("A444:D444").Copy() | Out-Null
$ws.Range("A445:D470").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new data rows (445:470 -> A=444..469, B/C/D from source data)
$ws.Range("A445").Value = 444
$ws.Range("B445").Value = 3435877
$ws.Range("C445").Value = 10279
$ws.Range("D445").Value = 84020
$ws.Range("A446").Value = 445
$ws.Range("B446").Value = 3448182
$ws.Range("C446").Value = 12305
$ws.Range("D446").Value = 84285
$ws.Range("A447").Value = 446
$ws.Range("B447").Value = 3469448
$ws.Range("C447").Value = 21266
$ws.Range("D447").Value = 84593
$ws.Range("A448").Value = 447
$ws.Range("B448").Value = 3486462
$ws.Range("C448").Value = 17014
$ws.Range("D448").Value = 84811
$ws.Range("A449").Value = 448
$ws.Range("B449").Value = 3504012
$ws.Range("C449").Value = 17550
$ws.Range("D449").Value = 85056
$ws.Range("A450").Value = 449
$ws.Range("B450").Value = 3519250
$ws.Range("C450").Value = 15238
$ws.Range("D450").Value = 85252
$ws.Range("A451").Value = 450
$ws.Range("B451").Value = 3527540
$ws.Range("C451").Value = 8290
$ws.Range("D451").Value = 85371
$ws.Range("A452").Value = 451
$ws.Range("B452").Value = 3535354
$ws.Range("C452").Value = 7814
$ws.Range("D452").Value = 85481
$ws.Range("A453").Value = 452
$ws.Range("B453").Value = 3544315
$ws.Range("C453").Value = 8961
$ws.Range("D453").Value = 85757
$ws.Range("A454").Value = 453
$ws.Range("B454").Value = 3558148
$ws.Range("C454").Value = 13833
$ws.Range("D454").Value = 86009
$ws.Range("A455").Value = 454
$ws.Range("B455").Value = 3575644
$ws.Range("C455").Value = 17496
$ws.Range("D455").Value = 86276
$ws.Range("A456").Value = 455
$ws.Range("B456").Value = 3585891
$ws.Range("C456").Value = 10247
$ws.Range("D456").Value = 86481
$ws.Range("A457").Value = 456
$ws.Range("B457").Value = 3595872
$ws.Range("C457").Value = 9981
$ws.Range("D457").Value = 86669
$ws.Range("A458").Value = 457
$ws.Range("B458").Value = 3602939
$ws.Range("C458").Value = 7067
$ws.Range("D458").Value = 86731
$ws.Range("A459").Value = 458
$ws.Range("B459").Value = 3608292
$ws.Range("C459").Value = 5353
$ws.Range("D459").Value = 86870
$ws.Range("A460").Value = 459
$ws.Range("B460").Value = 3615896
$ws.Range("C460").Value = 7604
$ws.Range("D460").Value = 87148
$ws.Range("A461").Value = 460
$ws.Range("B461").Value = 3627777
$ws.Range("C461").Value = 11881
$ws.Range("D461").Value = 87405
$ws.Range("A462").Value = 461
$ws.Range("B462").Value = 3638504
$ws.Range("C462").Value = 10727
$ws.Range("D462").Value = 87639
$ws.Range("A463").Value = 462
$ws.Range("B463").Value = 3646600
$ws.Range("C463").Value = 8096
$ws.Range("D463").Value = 87852
$ws.Range("A464").Value = 463
$ws.Range("B464").Value = 3653019
$ws.Range("C464").Value = 6419
$ws.Range("D464").Value = 87960
$ws.Range("A465").Value = 464
$ws.Range("B465").Value = 3654201
$ws.Range("C465").Value = 1182
$ws.Range("D465").Value = 87973
$ws.Range("A466").Value = 465
$ws.Range("B466").Value = 3659990
$ws.Range("C466").Value = 5789
$ws.Range("D466").Value = 88039
$ws.Range("A467").Value = 466
$ws.Range("B467").Value = 3662568
$ws.Range("C467").Value = 2578
$ws.Range("D467").Value = 88161
$ws.Range("A468").Value = 467
$ws.Range("B468").Value = 3667041
$ws.Range("C468").Value = 4473
$ws.Range("D468").Value = 88479
$ws.Range("A469").Value = 468
$ws.Range("B469").Value = 3673969
$ws.Range("C469").Value = 6928
$ws.Range("D469").Value = 88689
$ws.Range("A470").Value = 469
$ws.Range("B470").Value = 3680159
$ws.Range("C470").Value = 6190
$ws.Range("D470").Value = 88884

# Restore the view/selection state to match the post-edit workbook
$ws.Range("F431").Select()
